$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.636.67"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.44%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.295.63"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.39%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.62"
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.69%  "

$ws.Range("E7").Value = "  +0.78%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("E9").Value = "  +3.86%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.37"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +13.18%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0807"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.31%  "

$ws.Range("E12").Value = "  -1.38%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.73"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.53%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.645.61"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.47%  "

$ws.Range("E15").Value = "  +3.28%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.288.66"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.07%  "

$ws.Range("E17").Value = "  +5.77%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.491.22"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.31%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.74"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.21%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0920"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.03%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.03"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.14%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.96"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.94%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "243.40"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.29%  "

$ws.Range("E24").Value = "  +0.82%  "

$ws.Range("E25").Value = "  +2.42%  "

$ws.Range("E26").Value = "  -0.11%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.08"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.61%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "36.89"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.85%  "

$ws.Range("E29").Value = "  +0.95%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.12"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.35%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "161.09"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.79%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.34"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.86%  "

$ws.Range("E33").Value = "  +0.01%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0753"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.51%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.12"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.51%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.39"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.77%  "

$ws.Range("E37").Value = "  +3.21%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.89"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.14%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.38"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.39%  "

$ws.Range("E40").Value = "  -0.33%  "

$ws.Range("E41").Value = "  +6.35%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.45"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +20.03%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.008.44"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.89%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.41"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.54%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0287"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.36%  "

$ws.Range("E46").Value = "  +5.13%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.21"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.94%  "

$ws.Range("E48").Value = "  +4.71%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.56"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.34%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "72.83"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.21%  "

$ws.Range("E51").Value = "  -0.57%  "
